# Refresh crypto price/volume snapshot (D,E columns) for rows 2-51.
# Values are stored as literal text in the sheet; numeric-looking prices
# are prefixed with a text-quote marker ($q) so Excel keeps them as strings
# instead of silently converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$q = "'"

$ws.Range("D2").Value = '64.088.02'
$ws.Range("E2").Value = '  -3.28%  '

$ws.Range("D3").Value = '3.129.95'
$ws.Range("E3").Value = '  -2.25%  '

$ws.Range("D4").Value = $q + '1.00'
$ws.Range("E4").Value = '  +0.11%  '

$ws.Range("D5").Value = $q + '606.12'
$ws.Range("E5").Value = '  -0.26%  '

$ws.Range("D6").Value = $q + '148.47'
$ws.Range("E6").Value = '  -5.17%  '

$ws.Range("E7").Value = '  +0.08%  '

$ws.Range("D8").Value = '3.128.26'
$ws.Range("E8").Value = '  -2.23%  '

$ws.Range("D9").Value = $q + '0.531'
$ws.Range("E9").Value = '  -4.09%  '

$ws.Range("D10").Value = $q + '0.152'
$ws.Range("E10").Value = '  -5.11%  '

$ws.Range("D11").Value = $q + '5.58'
$ws.Range("E11").Value = '  -1.18%  '

$ws.Range("D12").Value = $q + '0.476'
$ws.Range("E12").Value = '  -5.38%  '

$ws.Range("D13").Value = $q + '0.0000257'
$ws.Range("E13").Value = '  -4.45%  '

$ws.Range("D14").Value = $q + '36.62'
$ws.Range("E14").Value = '  -4.66%  '

$ws.Range("D15").Value = '3.646.30'
$ws.Range("E15").Value = '  -2.26%  '

$ws.Range("D16").Value = '64.228.92'
$ws.Range("E16").Value = '  -3.25%  '

$ws.Range("E17").Value = '  +0.03%  '

$ws.Range("D18").Value = '3.137.71'
$ws.Range("E18").Value = '  -1.75%  '

$ws.Range("D19").Value = $q + '6.97'
$ws.Range("E19").Value = '  -4.72%  '

$ws.Range("D20").Value = $q + '479.78'
$ws.Range("E20").Value = '  -5.45%  '

$ws.Range("D21").Value = $q + '14.53'
$ws.Range("E21").Value = '  -4.99%  '

$ws.Range("D22").Value = $q + '0.708'
$ws.Range("E22").Value = '  -3.14%  '

$ws.Range("D23").Value = $q + '7.73'
$ws.Range("E23").Value = '  -3.29%  '

$ws.Range("D24").Value = $q + '13.63'
$ws.Range("E24").Value = '  -6.68%  '

$ws.Range("D25").Value = $q + '83.68'
$ws.Range("E25").Value = '  -1.68%  '

$ws.Range("E26").Value = '  +0.02%  '

$ws.Range("D27").Value = $q + '2.93'
$ws.Range("E27").Value = '  -2.48%  '

$ws.Range("D28").Value = $q + '8.51'
$ws.Range("E28").Value = '  -5.80%  '

$ws.Range("D29").Value = $q + '0.126'
$ws.Range("E29").Value = '  -1.60%  '

$ws.Range("D30").Value = $q + '2.23'
$ws.Range("E30").Value = '  -5.23%  '

$ws.Range("D31").Value = $q + '6.94'
$ws.Range("E31").Value = '  -0.74%  '

$ws.Range("D32").Value = $q + '1.00'
$ws.Range("E32").Value = '  -0.06%  '

$ws.Range("D33").Value = $q + '2.72'
$ws.Range("E33").Value = '  -6.95%  '

$ws.Range("D34").Value = $q + '26.65'
$ws.Range("E34").Value = '  -5.61%  '

$ws.Range("E35").Value = '  -5.16%  '

$ws.Range("D36").Value = $q + '6.08'
$ws.Range("E36").Value = '  -5.75%  '

$ws.Range("D37").Value = $q + '54.53'
$ws.Range("E37").Value = '  -1.56%  '

$ws.Range("D38").Value = $q + '3.23'
$ws.Range("E38").Value = '  +6.57%  '

$ws.Range("D39").Value = '0.0₃0744'
$ws.Range("E39").Value = '  -3.30%  '

$ws.Range("D40").Value = $q + '447.89'
$ws.Range("E40").Value = '  -10.51%  '

$ws.Range("D41").Value = $q + '0.0399'
$ws.Range("E41").Value = '  -5.20%  '

$ws.Range("E42").Value = '  -5.00%  '

$ws.Range("D43").Value = $q + '8.45'
$ws.Range("E43").Value = '  -3.06%  '

$ws.Range("D44").Value = '2.868.68'
$ws.Range("E44").Value = '  -1.30%  '

$ws.Range("D45").Value = $q + '0.271'
$ws.Range("E45").Value = '  -8.33%  '

$ws.Range("D46").Value = $q + '2.31'
$ws.Range("E46").Value = '  -4.90%  '

$ws.Range("D47").Value = $q + '26.52'
$ws.Range("E47").Value = '  -5.95%  '

$ws.Range("D48").Value = $q + '0.998'
$ws.Range("E48").Value = '  -0.03%  '

$ws.Range("D49").Value = $q + '0.115'
$ws.Range("E49").Value = '  -1.65%  '

$ws.Range("D50").Value = $q + '2.30'
$ws.Range("E50").Value = '  -3.73%  '

$ws.Range("D51").Value = $q + '118.93'
$ws.Range("E51").Value = '  -3.24%  '
